$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.523.72'
$ws.Range("D3").Value = '3.369.75'
$ws.Range("E3").Value = '  +4.66%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '191.77'
$ws.Range("E5").Value = '  +5.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '593.81'
$ws.Range("E6").Value = '  +2.74%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.607'
$ws.Range("E8").Value = '  +0.99%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.134'
$ws.Range("E9").Value = '  +3.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.75'
$ws.Range("E10").Value = '  +2.64%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.421'
$ws.Range("E11").Value = '  +2.68%  '
$ws.Range("D12").Value = '3.960.31'
$ws.Range("E12").Value = '  +4.80%  '
$ws.Range("E13").Value = '  +1.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.64'
$ws.Range("E14").Value = '  +3.63%  '
$ws.Range("D15").Value = '69.545.52'
$ws.Range("E15").Value = '  +3.09%  '
$ws.Range("E16").Value = '  +2.25%  '
$ws.Range("D17").Value = '3.360.29'
$ws.Range("E17").Value = '  +3.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '448.20'
$ws.Range("E18").Value = '  +13.73%  '
$ws.Range("E19").Value = '  +1.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.80'
$ws.Range("E20").Value = '  +3.20%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.81'
$ws.Range("E21").Value = '  +3.70%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.52'
$ws.Range("E22").Value = '  +3.85%  '
$ws.Range("E23").Value = '  +0.20%  '
$ws.Range("D24").Value = '3.519.74'
$ws.Range("E24").Value = '  +4.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.518'
$ws.Range("E25").Value = '  +1.16%  '
$ws.Range("E26").Value = '  +4.22%  '
$ws.Range("E27").Value = '  +4.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.59'
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  +2.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '23.24'
$ws.Range("E31").Value = '  +2.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.62'
$ws.Range("E32").Value = '  +1.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.30'
$ws.Range("E33").Value = '  +4.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.03'
$ws.Range("E34").Value = '  +1.17%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  +3.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '164.81'
$ws.Range("E37").Value = '  +2.31%  '
$ws.Range("E38").Value = '  +3.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.18'
$ws.Range("E39").Value = '  +3.84%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.822'
$ws.Range("E40").Value = '  +2.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.61'
$ws.Range("E41").Value = '  +1.34%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.50'
$ws.Range("E42").Value = '  +0.42%  '
$ws.Range("D43").Value = '2.743.01'
$ws.Range("E43").Value = '  +5.70%  '
$ws.Range("E44").Value = '  +2.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '25.54'
$ws.Range("E45").Value = '  +4.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0689'
$ws.Range("E46").Value = '  +0.99%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '345.20'
$ws.Range("E47").Value = '  +4.00%  '
$ws.Range("E48").Value = '  +0.80%  '
$ws.Range("E49").Value = '  +3.87%  '
$ws.Range("B50").Value = 'ONDO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.04'
$ws.Range("E50").Value = '  +8.39%  '
$ws.Range("B51").Value = 'Arweave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '32.87'
$ws.Range("E51").Value = '  +7.44%  '
